$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the datetime number format to the new "Bid end" data cells first, so
# it claims cellXfs style index 1 (matches the target style table order).
$ws.Range("D5:D7").NumberFormat = "yyyy-mm-dd h:mm:ss"

# New header cell for the added "Bid end" column.
$ws.Range("D1").Value = "Bid end"

# Row 2 - Item ID / Amount / Seconds (no Bid end value on this row).
# A leading apostrophe forces the numeric-looking IDs/values to be stored as
# text (matching the source data) instead of being auto-coerced to numbers;
# re-applying the "Normal" style afterwards clears the transient quote-prefix
# formatting that the apostrophe entry leaves behind.
$ws.Range("A2").Value = "'204498232262"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "£"
$ws.Range("C2").Value = "'5"
$ws.Range("C2").Style = "Normal"

# Row 3
$ws.Range("A3").Value = "'204498232262"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "£"
$ws.Range("C3").Value = "'5"
$ws.Range("C3").Style = "Normal"

# Row 4
$ws.Range("A4").Value = "'204498232262"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "£"
$ws.Range("C4").Value = "'5"
$ws.Range("C4").Style = "Normal"

# Row 5 - also carries a "Bid end" timestamp.
$ws.Range("A5").Value = "'204498232262"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "£"
$ws.Range("C5").Value = "'5"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = 45218.91527777778

# Row 6
$ws.Range("A6").Value = "'175971859943"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "£"
$ws.Range("C6").Value = "'5"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = 45218.91319444445

# Row 7
$ws.Range("A7").Value = "'175971859943"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "£"
$ws.Range("C7").Value = "'5"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = 45218.91527777778

# Move the active selection to H13, matching the saved view state.
$ws.Range("H13").Select()
